# Weekly price-data refresh for Fruta / hortaliza (Durazno, Agricola del Norte S.A. de Arica).
# The sheet row order is fixed (row 1 = headers); column D and K..T values for each data
# row (2-27) are updated in place to reflect this week's refreshed price/variety figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44236
$ws.Range("K2").Value = 'Doctor Davis'
$ws.Range("L2").Value = 'Segunda'
$ws.Range("N2").Value = 20000
$ws.Range("O2").Value = 21000
$ws.Range("P2").Value = 20500
$ws.Range("Q2").Value = '$/caja 18 kilos granel'
$ws.Range("S2").Value = 1139
# Row 3
$ws.Range("D3").Value = 44174
$ws.Range("K3").Value = 'Kurakata'
$ws.Range("N3").Value = 20000
$ws.Range("O3").Value = 21000
$ws.Range("P3").Value = 20500
$ws.Range("R3").Value = 'Región de Coquimbo'
$ws.Range("S3").Value = 1139
# Row 4
$ws.Range("D4").Value = 44258
$ws.Range("K4").Value = 'Doctor Davis'
$ws.Range("M4").Value = 300
$ws.Range("N4").Value = 9000
$ws.Range("O4").Value = 10000
$ws.Range("P4").Value = 9500
$ws.Range("Q4").Value = '$/bandeja 10 kilos empedrada'
$ws.Range("R4").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S4").Value = 950
$ws.Range("T4").Value = 10
# Row 5
$ws.Range("D5").Value = 44203
$ws.Range("K5").Value = 'Carson'
$ws.Range("L5").Value = 'Tercera'
$ws.Range("M5").Value = 270
$ws.Range("Q5").Value = '$/caja 18 kilos granel'
# Row 6
$ws.Range("D6").Value = 44203
$ws.Range("K6").Value = 'Flavor Crest'
$ws.Range("M6").Value = 250
$ws.Range("Q6").Value = '$/caja 18 kilos granel'
# Row 7
$ws.Range("D7").Value = 44257
$ws.Range("K7").Value = 'September Sweet'
$ws.Range("N7").Value = 19000
$ws.Range("O7").Value = 20000
$ws.Range("P7").Value = 19500
$ws.Range("Q7").Value = '$/caja 18 kilos granel'
$ws.Range("S7").Value = 1083
# Row 8
$ws.Range("D8").Value = 44223
$ws.Range("K8").Value = 'Andross'
$ws.Range("N8").Value = 20000
$ws.Range("O8").Value = 21000
$ws.Range("P8").Value = 20500
$ws.Range("Q8").Value = '$/caja 18 kilos granel'
$ws.Range("R8").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S8").Value = 1139
# Row 9
$ws.Range("D9").Value = 44209
$ws.Range("K9").Value = 'Carson'
$ws.Range("M9").Value = 300
# Row 10
$ws.Range("D10").Value = 44209
$ws.Range("L10").Value = 'Tercera'
$ws.Range("N10").Value = 18000
$ws.Range("O10").Value = 19000
$ws.Range("P10").Value = 18500
$ws.Range("Q10").Value = '$/caja 18 kilos granel'
$ws.Range("S10").Value = 1028
# Row 11
$ws.Range("D11").Value = 44167
$ws.Range("K11").Value = 'Florida King'
$ws.Range("M11").Value = 250
$ws.Range("N11").Value = 25000
$ws.Range("O11").Value = 26000
$ws.Range("P11").Value = 25500
$ws.Range("S11").Value = 1417
# Row 12
$ws.Range("D12").Value = 44210
$ws.Range("K12").Value = 'Carson'
$ws.Range("M12").Value = 300
$ws.Range("N12").Value = 19000
$ws.Range("O12").Value = 20000
$ws.Range("P12").Value = 19500
$ws.Range("S12").Value = 1083
# Row 13
$ws.Range("D13").Value = 44210
$ws.Range("K13").Value = 'Rich Lady'
$ws.Range("N13").Value = 19000
$ws.Range("O13").Value = 20000
$ws.Range("P13").Value = 19500
$ws.Range("Q13").Value = '$/bandeja 18 kilos granel'
$ws.Range("R13").Value = 'Región de O''Higgins'
$ws.Range("S13").Value = 1083
$ws.Range("T13").Value = 18
# Row 14
$ws.Range("D14").Value = 44210
$ws.Range("K14").Value = 'Royal Glory'
$ws.Range("M14").Value = 300
$ws.Range("N14").Value = 19000
$ws.Range("O14").Value = 20000
$ws.Range("P14").Value = 19500
$ws.Range("Q14").Value = '$/bandeja 18 kilos granel'
$ws.Range("R14").Value = 'Región de O''Higgins'
$ws.Range("S14").Value = 1083
# Row 15
$ws.Range("D15").Value = 44161
$ws.Range("K15").Value = 'Florida King'
$ws.Range("M15").Value = 250
$ws.Range("N15").Value = 24000
$ws.Range("O15").Value = 25000
$ws.Range("P15").Value = 24500
$ws.Range("Q15").Value = '$/caja 18 kilos granel'
$ws.Range("S15").Value = 1361
# Row 16
$ws.Range("D16").Value = 44160
$ws.Range("K16").Value = 'Early Majestic'
$ws.Range("L16").Value = 'Segunda'
$ws.Range("M16").Value = 250
$ws.Range("N16").Value = 24000
$ws.Range("O16").Value = 25000
$ws.Range("P16").Value = 24500
$ws.Range("R16").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S16").Value = 1361
# Row 17
$ws.Range("D17").Value = 44216
$ws.Range("K17").Value = 'Andross'
$ws.Range("M17").Value = 270
$ws.Range("N17").Value = 17000
$ws.Range("O17").Value = 18000
$ws.Range("P17").Value = 17500
$ws.Range("Q17").Value = '$/caja 16 kilos empedrada'
$ws.Range("S17").Value = 1094
$ws.Range("T17").Value = 16
# Row 18
$ws.Range("D18").Value = 44222
$ws.Range("K18").Value = 'Elegant Lady'
$ws.Range("M18").Value = 300
$ws.Range("N18").Value = 20000
$ws.Range("O18").Value = 21000
$ws.Range("P18").Value = 20500
$ws.Range("S18").Value = 1139
# Row 19
$ws.Range("D19").Value = 44278
$ws.Range("K19").Value = 'Phillips Cling'
$ws.Range("M19").Value = 250
$ws.Range("N19").Value = 23000
$ws.Range("O19").Value = 24000
$ws.Range("P19").Value = 23500
$ws.Range("Q19").Value = '$/bandeja 18 kilos granel'
$ws.Range("S19").Value = 1306
# Row 20
$ws.Range("D20").Value = 44208
$ws.Range("K20").Value = 'Rich Lady'
$ws.Range("L20").Value = 'Primera'
$ws.Range("N20").Value = 28000
$ws.Range("O20").Value = 29000
$ws.Range("P20").Value = 28500
$ws.Range("S20").Value = 1583
# Row 21
$ws.Range("D21").Value = 44229
$ws.Range("K21").Value = 'Doctor Davis'
$ws.Range("L21").Value = 'Primera'
$ws.Range("M21").Value = 320
$ws.Range("Q21").Value = '$/bandeja 18 kilos granel'
# Row 22
$ws.Range("D22").Value = 44201
$ws.Range("L22").Value = 'Segunda'
$ws.Range("N22").Value = 22000
$ws.Range("O22").Value = 23000
$ws.Range("P22").Value = 22500
$ws.Range("S22").Value = 1250
# Row 23
$ws.Range("D23").Value = 44215
$ws.Range("K23").Value = 'Andross'
$ws.Range("M23").Value = 300
$ws.Range("N23").Value = 19000
$ws.Range("O23").Value = 20000
$ws.Range("P23").Value = 19500
$ws.Range("Q23").Value = '$/bandeja 18 kilos granel'
$ws.Range("S23").Value = 1083
# Row 24
$ws.Range("D24").Value = 44217
$ws.Range("K24").Value = 'Flavor Crest'
$ws.Range("L24").Value = 'Segunda'
$ws.Range("M24").Value = 250
$ws.Range("Q24").Value = '$/caja 18 kilos empedrada'
# Row 25
$ws.Range("D25").Value = 44217
$ws.Range("K25").Value = 'Royal Glory'
$ws.Range("L25").Value = 'Tercera'
$ws.Range("M25").Value = 270
$ws.Range("N25").Value = 17000
$ws.Range("O25").Value = 18000
$ws.Range("P25").Value = 17500
$ws.Range("Q25").Value = '$/bandeja 18 kilos granel'
$ws.Range("R25").Value = 'Región de O''Higgins'
$ws.Range("S25").Value = 972
$ws.Range("T25").Value = 18
# Row 26
$ws.Range("D26").Value = 44273
$ws.Range("K26").Value = 'Doctor Davis'
$ws.Range("M26").Value = 280
$ws.Range("N26").Value = 22000
$ws.Range("O26").Value = 23000
$ws.Range("P26").Value = 22500
$ws.Range("Q26").Value = '$/bandeja 18 kilos granel'
$ws.Range("S26").Value = 1250
# Row 27
$ws.Range("D27").Value = 44273
$ws.Range("K27").Value = 'Phillips Cling'
$ws.Range("Q27").Value = '$/bandeja 18 kilos granel'
